$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.520.79"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "'3.004.48"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'596.12"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").Value = "'144.33"
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'3.005.00"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "'34.37"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "'3.498.28"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "'7.04"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'61.537.28"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").Value = "'3.000.84"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "'454.51"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "'14.06"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'0.688"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'7.37"
$ws.Range("D24").Value = "'81.97"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").Value = "'10.79"
$ws.Range("E26").Value = "  +6.05%  "
$ws.Range("D27").Value = "'11.99"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "'27.55"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").Value = "'0.0₃0839"
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").Value = "'5.79"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").Value = "'50.37"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  +9.87%  "
$ws.Range("D42").Value = "'2.91"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "'400.04"
$ws.Range("E43").Value = "  -5.81%  "
$ws.Range("D44").Value = "'39.91"
$ws.Range("E44").Value = "  +5.10%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.271"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0354"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'2.719.22"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").Value = "'132.30"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "'0.108"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").Value = "  +1.58%  "
